$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035673440647108
$ws.Range("D2").Value = 1.032822811053197
$ws.Range("E2").Value = 1.034619009003198
$ws.Range("F2").Value = 1.034346543335745
$ws.Range("I2").Value = 1.035267872505472
$ws.Range("J2").Value = 1.040785930075492
$ws.Range("K2").Value = 1.035626932582745
$ws.Range("L2").Value = 1.03741796361989
$ws.Range("M2").Value = 1.037146280468248
$ws.Range("N2").Value = 1.042263965476052

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03730602150839
$ws.Range("D3").Value = 1.033456399407444
$ws.Range("E3").Value = 1.036029250861922
$ws.Range("F3").Value = 1.036617500872505
$ws.Range("I3").Value = 1.035614217647267
$ws.Range("J3").Value = 1.042058713845839
$ws.Range("K3").Value = 1.036069841880408
$ws.Range("L3").Value = 1.038635831029
$ws.Range("M3").Value = 1.039222517231887
$ws.Range("N3").Value = 1.043538556745305

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038359123688883
$ws.Range("D4").Value = 1.033865498509385
$ws.Range("E4").Value = 1.0369390799284
$ws.Range("F4").Value = 1.038083082856138
$ws.Range("I4").Value = 1.03583606318397
$ws.Range("J4").Value = 1.042878708986815
$ws.Range("K4").Value = 1.036354873154433
$ws.Range("L4").Value = 1.039420672106835
$ws.Range("M4").Value = 1.040561791040016
$ws.Range("N4").Value = 1.044359716373436

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038801076634731
$ws.Range("D5").Value = 1.034037276700652
$ws.Range("E5").Value = 1.037320940114961
$ws.Range("F5").Value = 1.038698311826356
$ws.Range("I5").Value = 1.035928788763703
$ws.Range("J5").Value = 1.043222590330306
$ws.Range("K5").Value = 1.03647432941281
$ws.Range("L5").Value = 1.039749864726716
$ws.Range("M5").Value = 1.041123843470272
$ws.Range("N5").Value = 1.044704086067878

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038875237706621
$ws.Range("D6").Value = 1.034066106940558
$ws.Range("E6").Value = 1.037385019421646
$ws.Range("F6").Value = 1.038801559381876
$ws.Range("I6").Value = 1.035944326331172
$ws.Range("J6").Value = 1.043280280358805
$ws.Range("K6").Value = 1.036494364973751
$ws.Range("L6").Value = 1.039805093703949
$ws.Range("M6").Value = 1.041218157897136
$ws.Range("N6").Value = 1.044761858022835

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038365032094759
$ws.Range("D7").Value = 1.033867794628209
$ws.Range("E7").Value = 1.036944184829022
$ws.Range("F7").Value = 1.03809130708225
$ws.Range("I7").Value = 1.035837304298004
$ws.Range("J7").Value = 1.042883307247575
$ws.Range("K7").Value = 1.03635647078831
$ws.Range("L7").Value = 1.03942507374375
$ws.Range("M7").Value = 1.04056930502243
$ws.Range("N7").Value = 1.044364321164254

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036225869019525
$ws.Range("D8").Value = 1.03303711781638
$ws.Range("E8").Value = 1.035096171712946
$ws.Range("F8").Value = 1.03511484368188
$ws.Range("I8").Value = 1.035385392121637
$ws.Range("J8").Value = 1.041216822479269
$ws.Range("K8").Value = 1.035776940657547
$ws.Range("L8").Value = 1.037830217365077
$ws.Range("M8").Value = 1.037848837059198
$ws.Range("N8").Value = 1.042695469796428

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032430510189243
$ws.Range("D9").Value = 1.031566548796502
$ws.Range("E9").Value = 1.031818558493934
$ws.Range("F9").Value = 1.029839002309584
$ws.Range("I9").Value = 1.034571576662682
$ws.Range("J9").Value = 1.038252262617374
$ws.Range("K9").Value = 1.034743657214135
$ws.Range("L9").Value = 1.03499483727929
$ws.Range("M9").Value = 1.03302180991313
$ws.Range("N9").Value = 1.03972669991949

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02988187524849
$ws.Range("D10").Value = 1.030581433259586
$ws.Range("E10").Value = 1.029618451384772
$ws.Range("F10").Value = 1.02629926088867
$ws.Range("I10").Value = 1.034017052358593
$ws.Range("J10").Value = 1.036256254478592
$ws.Range("K10").Value = 1.034046504106044
$ws.Range("L10").Value = 1.033087006506328
$ws.Range("M10").Value = 1.029779880914249
$ws.Range("N10").Value = 1.037727857220205

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028773701601263
$ws.Range("D11").Value = 1.030153707707681
$ws.Range("E11").Value = 1.028662037278791
$ws.Range("F11").Value = 1.024760789213491
$ws.Range("I11").Value = 1.033774045258514
$ws.Range("J11").Value = 1.035387122067489
$ws.Range("K11").Value = 1.033742620617297
$ws.Range("L11").Value = 1.032256558765104
$ws.Range("M11").Value = 1.02837006466518
$ws.Range("N11").Value = 1.036857490541388

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028361366396811
$ws.Range("D12").Value = 1.029994653272127
$ws.Range("E12").Value = 1.028306203133305
$ws.Range("F12").Value = 1.024188436633395
$ws.Range("I12").Value = 1.033683342330666
$ws.Range("J12").Value = 1.035063543476395
$ws.Range("K12").Value = 1.033629438962836
$ws.Range("L12").Value = 1.031947426286469
$ws.Range("M12").Value = 1.027845458113516
$ws.Range("N12").Value = 1.036533452431579

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028449846148644
$ws.Range("D13").Value = 1.030028779101587
$ws.Range("E13").Value = 1.028382557233675
$ws.Range("F13").Value = 1.024311249276397
$ws.Range("I13").Value = 1.033702818366167
$ws.Range("J13").Value = 1.035132986111415
$ws.Range("K13").Value = 1.033653730707964
$ws.Range("L13").Value = 1.032013766680375
$ws.Range("M13").Value = 1.027958030943792
$ws.Range("N13").Value = 1.036602993683106

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028739632445009
$ws.Range("D14").Value = 1.030140563871035
$ws.Range("E14").Value = 1.028632635826496
$ws.Range("F14").Value = 1.024713496836946
$ws.Range("I14").Value = 1.033766556713
$ws.Range("J14").Value = 1.03536039023215
$ws.Range("K14").Value = 1.033733271245513
$ws.Range("L14").Value = 1.032231019473587
$ws.Range("M14").Value = 1.028326719838228
$ws.Range("N14").Value = 1.036830720743777

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028918084657301
$ws.Range("D15").Value = 1.030209414455941
$ws.Range("E15").Value = 1.02878664026307
$ws.Range("F15").Value = 1.024961215141987
$ws.Range("I15").Value = 1.033805769671111
$ws.Range("J15").Value = 1.035500402337151
$ws.Range("K15").Value = 1.033782238166749
$ws.Range("L15").Value = 1.032364787240052
$ws.Range("M15").Value = 1.028553756012755
$ws.Range("N15").Value = 1.036970931682026

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029955322562464
$ws.Range("D16").Value = 1.030609795196228
$ws.Range("E16").Value = 1.029681845046516
$ws.Range("F16").Value = 1.026401240415289
$ws.Range("I16").Value = 1.034033118581953
$ws.Range("J16").Value = 1.036313832454857
$ws.Range("K16").Value = 1.034066629166565
$ws.Range("L16").Value = 1.033142027821889
$ws.Range("M16").Value = 1.029873315922474
$ws.Range("N16").Value = 1.037785516963801

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030604709904113
$ws.Range("D17").Value = 1.030860629549275
$ws.Range("E17").Value = 1.030242367651369
$ws.Range("F17").Value = 1.027302970036378
$ws.Range("I17").Value = 1.034174950532276
$ws.Range("J17").Value = 1.036822766720583
$ws.Range("K17").Value = 1.034244479061495
$ws.Range("L17").Value = 1.033628396997738
$ws.Range("M17").Value = 1.030699402522268
$ws.Range("N17").Value = 1.038295173974561

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030983044246052
$ws.Range("D18").Value = 1.031006824895426
$ws.Range("E18").Value = 1.030568949930325
$ws.Range("F18").Value = 1.027828382085686
$ws.Range("I18").Value = 1.034257399742129
$ws.Range("J18").Value = 1.037119152811708
$ws.Range("K18").Value = 1.034348022052124
$ws.Range("L18").Value = 1.033911669725752
$ws.Range("M18").Value = 1.031180664216846
$ws.Range("N18").Value = 1.03859198096793

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031111972039194
$ws.Range("D19").Value = 1.031056654815027
$ws.Range("E19").Value = 1.03068024526872
$ws.Range("F19").Value = 1.028007441459057
$ws.Range("I19").Value = 1.034285465606877
$ws.Range("J19").Value = 1.037220134131292
$ws.Range("K19").Value = 1.034383294763914
$ws.Range("L19").Value = 1.034008187910312
$ws.Range("M19").Value = 1.031344664451607
$ws.Range("N19").Value = 1.038693105692571

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030535082677793
$ws.Range("D20").Value = 1.030833729014243
$ws.Range("E20").Value = 1.030182266339712
$ws.Range("F20").Value = 1.027206280274378
$ws.Range("I20").Value = 1.034159762193019
$ws.Range("J20").Value = 1.036768211272628
$ws.Range("K20").Value = 1.034225417526922
$ws.Range("L20").Value = 1.033576257547392
$ws.Range("M20").Value = 1.030610831527383
$ws.Range("N20").Value = 1.038240541051611

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028654317364658
$ws.Range("D21").Value = 1.030107650985976
$ws.Range("E21").Value = 1.02855901003262
$ws.Range("F21").Value = 1.024595069951474
$ws.Range("I21").Value = 1.033747799520339
$ws.Range("J21").Value = 1.035293446055719
$ws.Range("K21").Value = 1.033709857010226
$ws.Range("L21").Value = 1.03216706247673
$ws.Range("M21").Value = 1.028218176290458
$ws.Range("N21").Value = 1.036763681498936

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027467686941228
$ws.Range("D22").Value = 1.029650104701459
$ws.Range("E22").Value = 1.027535044526235
$ws.Range("F22").Value = 1.022948099308599
$ws.Range("I22").Value = 1.033486238512773
$ws.Range("J22").Value = 1.03436188941887
$ws.Range("K22").Value = 1.033383932588302
$ws.Range("L22").Value = 1.03127717798619
$ws.Range("M22").Value = 1.02670837636156
$ws.Range("N22").Value = 1.035830801944812

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028097139120164
$ws.Range("D23").Value = 1.029892757585326
$ws.Range("E23").Value = 1.028078191715311
$ws.Range("F23").Value = 1.023821693657119
$ws.Range("I23").Value = 1.033625139560014
$ws.Range("J23").Value = 1.034856139393507
$ws.Range("K23").Value = 1.033556880378357
$ws.Range("L23").Value = 1.031749294095095
$ws.Range("M23").Value = 1.027509276111544
$ws.Range("N23").Value = 1.036325753811104

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030566545590606
$ws.Range("D24").Value = 1.030845884554605
$ws.Range("E24").Value = 1.030209424650861
$ws.Range("F24").Value = 1.027249971919159
$ws.Range("I24").Value = 1.034166626011708
$ws.Range("J24").Value = 1.036792863972868
$ws.Range("K24").Value = 1.034234031212918
$ws.Range("L24").Value = 1.033599818412583
$ws.Range("M24").Value = 1.03065085474067
$ws.Range("N24").Value = 1.038265228761514

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033414870039738
$ws.Range("D25").Value = 1.031947549169766
$ws.Range("E25").Value = 1.032668490741117
$ws.Range("F25").Value = 1.031206774766572
$ws.Range("I25").Value = 1.034784062626719
$ws.Range("J25").Value = 1.039022073233669
$ws.Range("K25").Value = 1.035012235035127
$ws.Range("L25").Value = 1.035730894099317
$ws.Range("M25").Value = 1.034273809566011
$ws.Range("N25").Value = 1.040497603755159
